$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Admin" column (G) for each user row: default "NO", with "YES"
# for rows 5 and 9 (Sophia Cathie / Emily Smith), replacing the old "Sowmya"
# placeholder value.
$adminValues = @{
    2 = "NO"
    3 = "NO"
    4 = "NO"
    5 = "YES"
    6 = "NO"
    7 = "NO"
    8 = "NO"
    9 = "YES"
    10 = "NO"
    11 = "NO"
    12 = "NO"
    13 = "NO"
    14 = "NO"
    15 = "NO"
}

foreach ($row in $adminValues.Keys) {
    $ws.Range("G$row").Value = $adminValues[$row]
}

$ws.Range("G15").Select()
